$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @{
  1  = @('93÷4=', '98÷7=', '97÷9=', '40÷9=', '68÷2=')
  5  = @('55÷4=', '85÷9=', '38÷6=', '20÷2=', '81÷4=')
  9  = @('49÷3=', '31÷5=', '28÷2=', '41÷3=', '32÷4=')
  13 = @('69÷2=', '84÷3=', '51÷5=', '94÷2=', '93÷9=')
  17 = @('73÷2=', '67÷4=', '24÷6=', '28÷6=', '42÷3=')
}

foreach ($r in $newValues.Keys) {
  $vals = $newValues[$r]
  for ($c = 1; $c -le $vals.Count; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $vals[$c - 1]
  }
}
